$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.307.01"
$ws.Range("E2").Value = "  -2.76%  "

$ws.Range("D3").Value = "2.413.17"
$ws.Range("E3").Value = "  -3.80%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.65"
$ws.Range("E5").Value = "  -4.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.94"
$ws.Range("E6").Value = "  -4.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("E8").Value = "  -3.51%  "

$ws.Range("D9").Value = "2.415.10"
$ws.Range("E9").Value = "  -3.87%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0954"
$ws.Range("E10").Value = "  -5.69%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  -4.49%  "

$ws.Range("E13").Value = "  -4.15%  "

$ws.Range("D14").Value = "2.840.70"
$ws.Range("E14").Value = "  -3.81%  "

$ws.Range("D15").Value = "57.246.75"
$ws.Range("E15").Value = "  -2.70%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.35"
$ws.Range("E16").Value = "  -5.93%  "

$ws.Range("E17").Value = "  -4.66%  "

$ws.Range("D18").Value = "2.415.54"
$ws.Range("E18").Value = "  -3.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.32"
$ws.Range("E19").Value = "  -6.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "313.10"
$ws.Range("E20").Value = "  -2.79%  "

$ws.Range("E21").Value = "  -4.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.61"
$ws.Range("E23").Value = "  -5.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.42"
$ws.Range("E24").Value = "  -2.50%  "

$ws.Range("E25").Value = "  -4.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("E27").Value = "  -2.57%  "

$ws.Range("E28").Value = "  -5.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.05"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("D30").Value = "0.0₃0716"
$ws.Range("E30").Value = "  -5.92%  "

$ws.Range("E31").Value = "  -5.37%  "

$ws.Range("E32").Value = "  -5.32%  "

$ws.Range("E33").Value = "  +1.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.65"
$ws.Range("E36").Value = "  -3.98%  "

$ws.Range("E37").Value = "  -7.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.84"
$ws.Range("E38").Value = "  -4.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.25"
$ws.Range("E39").Value = "  -2.17%  "

$ws.Range("E40").Value = "  -5.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.763"
$ws.Range("E41").Value = "  -4.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.35"
$ws.Range("E42").Value = "  -6.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "263.49"
$ws.Range("E43").Value = "  -6.14%  "

$ws.Range("E44").Value = "  -2.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.579"
$ws.Range("E45").Value = "  -4.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "121.69"
$ws.Range("E46").Value = "  -6.23%  "

$ws.Range("E47").Value = "  -2.96%  "

$ws.Range("E48").Value = "  -4.19%  "

$ws.Range("E49").Value = "  -4.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.33"
$ws.Range("E50").Value = "  -5.12%  "

$ws.Range("D51").Value = "1.683.74"
$ws.Range("E51").Value = "  -4.29%  "
